# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Map of worksheet name -> list of (row, newValue) pairs to apply to column F.
$updates = @{
    "展览" = @(
        @{ Row = 2;  Value = 142 },
        @{ Row = 3;  Value = 342 },
        @{ Row = 4;  Value = 442 },
        @{ Row = 5;  Value = 1746 },
        @{ Row = 7;  Value = 2204 },
        @{ Row = 11; Value = 4985 },
        @{ Row = 12; Value = 14 },
        @{ Row = 14; Value = 309 },
        @{ Row = 15; Value = 229 },
        @{ Row = 16; Value = 31 },
        @{ Row = 17; Value = 190 },
        @{ Row = 21; Value = 3951 },
        @{ Row = 22; Value = 718 },
        @{ Row = 23; Value = 692 },
        @{ Row = 26; Value = 110 },
        @{ Row = 27; Value = 122 },
        @{ Row = 29; Value = 13 },
        @{ Row = 31; Value = 584 },
        @{ Row = 32; Value = 11 },
        @{ Row = 34; Value = 1007 },
        @{ Row = 36; Value = 2539 },
        @{ Row = 38; Value = 14 }
    )
    "全部类型" = @(
        @{ Row = 2;  Value = 142 },
        @{ Row = 3;  Value = 342 },
        @{ Row = 4;  Value = 442 },
        @{ Row = 5;  Value = 1746 },
        @{ Row = 7;  Value = 2204 },
        @{ Row = 11; Value = 4985 },
        @{ Row = 12; Value = 14 },
        @{ Row = 14; Value = 309 },
        @{ Row = 15; Value = 229 },
        @{ Row = 16; Value = 31 },
        @{ Row = 17; Value = 190 },
        @{ Row = 21; Value = 3951 },
        @{ Row = 22; Value = 718 },
        @{ Row = 23; Value = 692 },
        @{ Row = 26; Value = 110 },
        @{ Row = 27; Value = 122 },
        @{ Row = 29; Value = 13 },
        @{ Row = 31; Value = 584 },
        @{ Row = 32; Value = 11 },
        @{ Row = 35; Value = 1007 },
        @{ Row = 37; Value = 2539 },
        @{ Row = 39; Value = 14 }
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($entry in $updates[$sheetName]) {
        $ws.Cells.Item($entry.Row, 6).Value = $entry.Value
    }
}
